# Updated cryptos list on Mon Feb 12 22:13:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a "Price" (column D) cell while preventing
# Excel's COM layer from auto-coercing numeric-looking text ("328.18", etc.)
# into a floating point Double. A leading apostrophe forces text entry, and
# resetting the Style back to "Normal" afterwards strips the transient
# "quote prefix" cell style so the cell's XML stays style-free, exactly as
# it was before the edit.
function Set-PriceText($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText $ws.Range("D2") "49.866.66"
$ws.Range("E2").Value = "  +3.70%  "

# Row 3 - Ethereum
Set-PriceText $ws.Range("D3") "2.637.83"
$ws.Range("E3").Value = "  +5.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-PriceText $ws.Range("D5") "328.18"
$ws.Range("E5").Value = "  +2.34%  "

# Row 6 - Solana
Set-PriceText $ws.Range("D6") "110.50"
$ws.Range("E6").Value = "  +2.93%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.48%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +3.79%  "

# Row 10 - Avalanche
Set-PriceText $ws.Range("D10") "40.56"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11 - Chainlink
$ws.Range("E11").Value = "  +2.10%  "

# Row 12 - Dogecoin
Set-PriceText $ws.Range("D12") "0.0820"
$ws.Range("E12").Value = "  +1.06%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.66%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.60%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText $ws.Range("D15") "3.048.87"
$ws.Range("E15").Value = "  +5.34%  "

# Row 16 - WrappedEther
Set-PriceText $ws.Range("D16") "2.623.83"
$ws.Range("E16").Value = "  +4.81%  "

# Row 17 - Polygon
Set-PriceText $ws.Range("D17") "0.880"
$ws.Range("E17").Value = "  +5.32%  "

# Row 18 - WrappedBTC
Set-PriceText $ws.Range("D18") "49.820.25"
$ws.Range("E18").Value = "  +3.91%  "

# Row 19 - ImmutableX
Set-PriceText $ws.Range("D19") "3.05"
$ws.Range("E19").Value = "  +10.75%  "

# Row 20 - InternetComputer(DFINITY)
Set-PriceText $ws.Range("D20") "13.33"
$ws.Range("E20").Value = "  +3.21%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.46%  "

# Row 22 - ShibaInu
Set-PriceText $ws.Range("D22") "0.0₃0960"
$ws.Range("E22").Value = "  +2.36%  "

# Row 23 - BitcoinCash
Set-PriceText $ws.Range("D23") "281.63"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24 - Litecoin
Set-PriceText $ws.Range("D24") "73.02"
$ws.Range("E24").Value = "  +2.19%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.75%  "

# Row 26 - EthereumClassic
Set-PriceText $ws.Range("D26") "26.72"
$ws.Range("E26").Value = "  +3.34%  "

# Row 27 - Dai
Set-PriceText $ws.Range("D27") "0.999"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +6.90%  "

# Rows 29/30 swapped: Cosmos <-> InjectiveProtocol
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceText $ws.Range("D29") "36.50"
$ws.Range("E29").Value = "  +3.59%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-PriceText $ws.Range("D30") "9.92"
$ws.Range("E30").Value = "  +2.74%  "

# Row 31 - Kaspa
Set-PriceText $ws.Range("D31") "0.143"
$ws.Range("E31").Value = "  +3.02%  "

# Row 32 - OKB
Set-PriceText $ws.Range("D32") "49.72"
$ws.Range("E32").Value = "  +0.42%  "

# Rows 33/34 swapped: Celestia <-> Filecoin
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText $ws.Range("D33") "5.43"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceText $ws.Range("D34") "19.48"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  -0.04%  "

# Row 36 - Hedera
Set-PriceText $ws.Range("D36") "0.0793"
$ws.Range("E36").Value = "  +1.54%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +6.23%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +2.31%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +8.16%  "

# Rows 40/41 swapped: Stellar <-> Monero
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-PriceText $ws.Range("D40") "123.88"
$ws.Range("E40").Value = "  +2.33%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText $ws.Range("D41") "0.113"
$ws.Range("E41").Value = "  +1.41%  "

# Rows 42/43 swapped: WEMIXToken <-> EnergySwap
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText $ws.Range("D42") "22.28"
$ws.Range("E42").Value = "  +4.92%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-PriceText $ws.Range("D43") "2.23"
$ws.Range("E43").Value = "  +0.71%  "

# Row 44 - VeChain
Set-PriceText $ws.Range("D44") "0.0313"
$ws.Range("E44").Value = "  +4.24%  "

# Row 45 - NEARProtocol
$ws.Range("E45").Value = "  +6.46%  "

# Row 46 - Maker
Set-PriceText $ws.Range("D46") "2.064.26"
$ws.Range("E46").Value = "  +2.53%  "

# Row 47 - ApeXProtocol
Set-PriceText $ws.Range("D47") "2.27"
$ws.Range("E47").Value = "  +13.60%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +7.86%  "

# Row 49 - FraxShare
Set-PriceText $ws.Range("D49") "9.05"
$ws.Range("E49").Value = "  +0.65%  "

# Row 51 - BitcoinSV
Set-PriceText $ws.Range("D51") "81.44"
$ws.Range("E51").Value = "  +1.68%  "
